# Update cryptos list data (auto-generated edit replicating upstream diff)
# Commit message: Updated cryptos list on Mon Nov  6 19:19:00 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $c = $ws.Range($cellRef)
    # Force the cell to text format so Excel does not reinterpret
    # numeric-looking strings (e.g. "2.00", "17.10") as numbers and
    # strip significant trailing zeros / formatting.
    $c.NumberFormat = "@"
    $c.Value = $value
    # Reset the cell style back to Normal/default so we don't leave a
    # stray style index behind (keeps output identical in style terms).
    $c.Style = "Normal"
}

Set-TextCell "D2" '35.297.48'
Set-TextCell "E2" '  -0.19%  '
Set-TextCell "D3" '1.910.68'
Set-TextCell "E4" '  +0.13%  '
Set-TextCell "E5" '  +9.17%  '
Set-TextCell "D6" '254.63'
Set-TextCell "E6" '  +3.37%  '
Set-TextCell "E7" '  +0.21%  '
Set-TextCell "D8" '40.56'
Set-TextCell "E8" '  -1.93%  '
Set-TextCell "D9" '0.367'
Set-TextCell "E9" '  +4.72%  '
Set-TextCell "D10" '52.68'
Set-TextCell "E10" '  -0.18%  '
Set-TextCell "D11" '0.0767'
Set-TextCell "E11" '  +6.51%  '
Set-TextCell "E12" '  -0.56%  '
Set-TextCell "D13" '2.186.52'
Set-TextCell "E13" '  -0.04%  '
Set-TextCell "D14" '12.76'
Set-TextCell "E14" '  +5.23%  '
Set-TextCell "D15" '0.720'
Set-TextCell "E15" '  +2.60%  '
Set-TextCell "B16" 'Polkadot'
Set-TextCell "C16" 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell "D16" '4.93'
Set-TextCell "E16" '  +1.17%  '
Set-TextCell "B17" 'WrappedEther'
Set-TextCell "C17" 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell "D17" '1.921.82'
Set-TextCell "E17" '  +0.66%  '
Set-TextCell "D18" '35.281.29'
Set-TextCell "E18" '  -0.21%  '
Set-TextCell "D19" '74.41'
Set-TextCell "E19" '  +2.55%  '
Set-TextCell "D20" '0.0₃0852'
Set-TextCell "E20" '  +3.18%  '
Set-TextCell "D21" '243.67'
Set-TextCell "E21" '  +1.56%  '
Set-TextCell "D22" '13.11'
Set-TextCell "E22" '  +4.54%  '
Set-TextCell "E23" '  +5.26%  '
Set-TextCell "E24" '  +0.18%  '
Set-TextCell "E25" '  +4.48%  '
Set-TextCell "D26" '2.39'
Set-TextCell "E26" '  +3.92%  '
Set-TextCell "D27" '167.15'
Set-TextCell "E27" '  -1.47%  '
Set-TextCell "E28" '  +2.53%  '
Set-TextCell "D29" '18.76'
Set-TextCell "E29" '  +1.46%  '
Set-TextCell "E30" '  +4.81%  '
Set-TextCell "D31" '4.128.87'
Set-TextCell "E31" '  +19.46%  '
Set-TextCell "E32" '  +4.80%  '
Set-TextCell "D33" '2.00'
Set-TextCell "E33" '  +14.33%  '
Set-TextCell "D34" '1.65'
Set-TextCell "E34" '  +23.45%  '
Set-TextCell "E35" '  +3.27%  '
Set-TextCell "E36" '  +2.31%  '
Set-TextCell "E37" '  +0.15%  '
Set-TextCell "D38" '0.911'
Set-TextCell "E38" '  -3.38%  '
Set-TextCell "E39" '  -0.59%  '
Set-TextCell "B40" 'InjectiveProtocol'
Set-TextCell "C40" 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextCell "D40" '17.10'
Set-TextCell "E40" '  +5.23%  '
Set-TextCell "B41" 'VeChain'
Set-TextCell "C41" 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell "D41" '0.0217'
Set-TextCell "E41" '  +4.16%  '
Set-TextCell "D42" '96.87'
Set-TextCell "E42" '  +7.37%  '
Set-TextCell "E43" '  +0.47%  '
Set-TextCell "D44" '0.0649'
Set-TextCell "E44" '  +1.43%  '
Set-TextCell "D45" '1.336.88'
Set-TextCell "E45" '  -0.26%  '
Set-TextCell "E46" '  +0.72%  '
Set-TextCell "E47" '  +1.03%  '
Set-TextCell "E48" '  +2.77%  '
Set-TextCell "E49" '  -0.83%  '
Set-TextCell "D50" '45.41'
Set-TextCell "E50" '  -5.40%  '
Set-TextCell "D51" '11.85'
Set-TextCell "E51" '  +14.39%  '
